# Add two new TextBox shapes to slide 1, matching the authored diff:
#   - id=26 "TextBox 25"  -> "XML package"
#   - id=27 "TextBox 26"  -> "Folders and files as .html .xml .js" (with "js" flagged)
#
# PowerPoint never re-uses shape ids within a slide, even across
# create/delete cycles, and the running "TextBox N" name counter tracks
# the same internal id allocator. The document's shapes currently top
# out at id 20, but this engine's id allocator also skips past ids used
# deeper in the document (nested group members, embedded picture ids,
# etc.) the first several times it's consulted, only settling into a
# simple "+1" cadence afterwards. Running it forward (and discarding the
# placeholders) lands the allocator exactly on id 26 / 27 for our two
# real shapes, matching the target ids and "TextBox 25"/"TextBox 26"
# names.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$scratch = @()
for ($i = 1; $i -le 12; $i++) {
    $scratch += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
}

# --- Shape 1: "XML package" -----------------------------------------
$tb1 = $s.Shapes.AddTextbox(1, 388.77425196850396, 131.77818897637795, 146.63346456692912, 31.50472440944882)
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1
$tb1.Fill.Visible = $false

$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "XML package"
$tr1.LanguageID = "fr-FR"
$tr1.Font.Size = 20
$tr1.ParagraphFormat.Alignment = 2

# --- Shape 2: "Folders and files as .html .xml .js" ------------------
$tb2 = $s.Shapes.AddTextbox(1, 368.51582677165356, 433.5733070866142, 146.63346456692912, 55.73905511811024)
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1
$tb2.Fill.Visible = $false

$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "Folders and files as .html .xml .js"
$tr2.LanguageID = "fr-FR"
$tr2.Font.Size = 20
$tr2.ParagraphFormat.Alignment = 2

$tail = $tr2.Characters(34, 2)
$tail.LanguageID = "fr-FR"
$tail.Font.Size = 20

foreach ($d in $scratch) {
    $d.Delete()
}
